# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.944.28"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "3.363.90"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.74"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.360.68"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.469"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -4.25%  "
$ws.Range("D13").Value = "3.940.51"
$ws.Range("E13").Value = "  -1.38%  "
$ws.Range("E14").Value = "  +1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.01"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("E16").Value = "  -5.05%  "
$ws.Range("D17").Value = "3.364.81"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "61.113.35"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.95"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "376.51"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.02%  "
$ws.Range("E23").Value = "  -3.80%  "
$ws.Range("D24").Value = "3.499.37"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -4.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.99"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.75"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.165"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.32%  "
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("E36").Value = "  -6.32%  "
$ws.Range("E37").Value = "  -3.47%  "
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0750"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.82%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("E45").Value = "  -2.31%  "
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.85"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.19"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").Value = "2.349.72"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.35"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.56%  "
